$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Menkomdigi Resmikan 1.194 Titik Kampung Internet, Wagub Sumut: Selaras dengan PHTC'
$ws.Range("B2").Value = 'Nusantara | Selasa, 30 September 2025 13:47 WIB Menkomdigi Resmikan 1.194 Titik Kampung Internet, Wagub Sumut: Selaras dengan PHTC'
$ws.Range("D2").Value = 'https://news.okezone.com/read/2025/09/30/340/3173459/menkomdigi-resmikan-1-194-titik-kampung-internet-wagub-sumut-selaras-dengan-phtc'

$ws.Range("A3").Value = 'Bright hingga Siti Nurhaliza Bakal Meriahkan Indonesian Television Awards 2025'
$ws.Range("B3").Value = 'TV Scoop | Selasa, 30 September 2025 13:30 WIB Bright hingga Siti Nurhaliza Bakal Meriahkan Indonesian Television Awards 2025'
$ws.Range("D3").Value = 'https://celebrity.okezone.com/read/2025/09/30/598/3173403/bright-hingga-siti-nurhaliza-bakal-meriahkan-indonesian-television-awards-2025'

$ws.Range("A4").Value = 'Orange Bond Oversubscribed, PNM Raih Dana Rp16 Triliun'
$ws.Range("B4").Value = 'Hot Issue | Selasa, 30 September 2025 13:27 WIB Orange Bond Oversubscribed, PNM Raih Dana Rp16 Triliun'
$ws.Range("D4").Value = 'https://economy.okezone.com/read/2025/09/30/320/3173458/orange-bond-oversubscribed-pnm-raih-dana-rp16-triliun'

$ws.Range("A5").Value = 'Purbaya Alokasikan Rp479 Triliun untuk Subsidi Energi dan Kompensasi 2025'
$ws.Range("B5").Value = 'Hot Issue | Selasa, 30 September 2025 13:10 WIB Purbaya Alokasikan Rp479 Triliun untuk Subsidi Energi dan Kompensasi 2025'
$ws.Range("D5").Value = 'https://economy.okezone.com/read/2025/09/30/320/3173457/purbaya-alokasikan-rp479-triliun-untuk-subsidi-energi-dan-kompensasi-2025'

$ws.Range("A6").Value = 'Patrick Kluivert Bahagia, Mathew Baker Resmi Dikontrak Melbourne City hingga 2028!'
$ws.Range("B6").Value = 'Sepakbola Dunia | Selasa, 30 September 2025 13:03 WIB Patrick Kluivert Bahagia, Mathew Baker Resmi Dikontrak Melbourne City hingga 2028!'
$ws.Range("D6").Value = 'https://bola.okezone.com/read/2025/09/30/51/3173455/patrick-kluivert-bahagia-mathew-baker-resmi-dikontrak-melbourne-city-hingga-2028'

$ws.Range("A7").Value = 'Kondisi Terkini Vika Kolesnaya Usai Melahirkan, Masih Dirawat dan Pemulihan'
$ws.Range("B7").Value = 'Hot Gossip | Selasa, 30 September 2025 13:02 WIB Kondisi Terkini Vika Kolesnaya Usai Melahirkan, Masih Dirawat dan Pemulihan'
$ws.Range("D7").Value = 'https://celebrity.okezone.com/read/2025/09/30/33/3173418/kondisi-terkini-vika-kolesnaya-usai-melahirkan-masih-dirawat-dan-pemulihan'

$ws.Range("A8").Value = 'Ini Cara dan Download Cheat Harvest Moon: Back to Nature untuk Android'
$ws.Range("B8").Value = 'Techno | Selasa, 30 September 2025 12:52 WIB Ini Cara dan Download Cheat Harvest Moon: Back to Nature untuk Android'
$ws.Range("D8").Value = 'https://ototekno.okezone.com/read/2025/09/30/16/3173454/ini-cara-dan-download-cheat-harvest-moon-back-to-nature-untuk-android'

$ws.Range("A9").Value = 'Warga NTB Dapat Diskon Tiket untuk MotoGP Mandalika 2025, 2 Kategori Termahal Sudah Ludes!'
$ws.Range("B9").Value = 'MotoGP | Selasa, 30 September 2025 12:44 WIB Warga NTB Dapat Diskon Tiket untuk MotoGP Mandalika 2025, 2 Kategori Termahal Sudah Ludes!'
$ws.Range("D9").Value = 'https://sports.okezone.com/read/2025/09/30/38/3173452/warga-ntb-dapat-diskon-tiket-untuk-motogp-mandalika-2025-2-kategori-termahal-sudah-ludes'

$ws.Range("A10").Value = 'Tampil Bak Pejabat, Bedu Akui Pinjam Kemeja Sang Kakak saat Hadir di PA Jaksel'
$ws.Range("B10").Value = 'Hot Gossip | Selasa, 30 September 2025 12:41 WIB Tampil Bak Pejabat, Bedu Akui Pinjam Kemeja Sang Kakak saat Hadir di PA Jaksel'
$ws.Range("D10").Value = 'https://celebrity.okezone.com/read/2025/09/30/33/3173451/tampil-bak-pejabat-bedu-akui-pinjam-kemeja-sang-kakak-saat-hadir-di-pa-jaksel'

$ws.Range("A11").Value = 'Live di RCTI! Ini Jadwal Siaran Langsung Persib Bandung vs Bangkok United di AFC Champions League 2 2025-2026'
$ws.Range("B11").Value = 'Liga Champion | Selasa, 30 September 2025 12:34 WIB Live di RCTI! Ini Jadwal Siaran Langsung Persib Bandung vs Bangkok United di AFC Champions League 2 2025-2026'
$ws.Range("D11").Value = 'https://bola.okezone.com/read/2025/09/30/261/3173450/live-di-rcti-ini-jadwal-siaran-langsung-persib-bandung-vs-bangkok-united-di-afc-champions-league-2-2025-2026'

$ws.Range("A12").Value = 'Azizah Salsha Pilih Main Padel saat Pratama Arhan Ucap Ikrar Talak'
$ws.Range("B12").Value = 'Hot Gossip | Selasa, 30 September 2025 12:30 WIB Azizah Salsha Pilih Main Padel saat Pratama Arhan Ucap Ikrar Talak'
$ws.Range("D12").Value = 'https://celebrity.okezone.com/read/2025/09/30/33/3173396/azizah-salsha-pilih-main-padel-saat-pratama-arhan-ucap-ikrar-talak'

$ws.Range("A13").Value = 'IHSG Sesi I Turun ke Level 8.096'
$ws.Range("B13").Value = 'Market Update | Selasa, 30 September 2025 12:28 WIB IHSG Sesi I Turun ke Level 8.096'
$ws.Range("D13").Value = 'https://economy.okezone.com/read/2025/09/30/278/3173449/ihsg-sesi-i-turun-ke-level-8-096'

$ws.Range("A14").Value = 'China Vonis Mati 11 Anggota Mafia Penguasa Pusat Penipuan Myanmar'
$ws.Range("B14").Value = 'International | Selasa, 30 September 2025 12:25 WIB China Vonis Mati 11 Anggota Mafia Penguasa Pusat Penipuan Myanmar'
$ws.Range("D14").Value = 'https://news.okezone.com/read/2025/09/30/18/3173448/china-vonis-mati-11-anggota-mafia-penguasa-pusat-penipuan-myanmar'

$ws.Range("A15").Value = 'Besaran Insentif Mobil Diusulkan Berdasarkan TKDN'
$ws.Range("B15").Value = 'Autos | Selasa, 30 September 2025 12:18 WIB Besaran Insentif Mobil Diusulkan Berdasarkan TKDN'
$ws.Range("D15").Value = 'https://ototekno.okezone.com/read/2025/09/30/15/3173447/besaran-insentif-mobil-diusulkan-berdasarkan-tkdn'

$ws.Range("A16").Value = 'Dahnil Anzar Sebut Kuota Haji Indonesia 2026 Capai 221 Ribu'
$ws.Range("B16").Value = 'Haji & Umroh | Selasa, 30 September 2025 12:14 WIB Dahnil Anzar Sebut Kuota Haji Indonesia 2026 Capai 221 Ribu'
$ws.Range("D16").Value = 'https://muslim.okezone.com/read/2025/09/30/398/3173446/dahnil-anzar-sebut-kuota-haji-indonesia-2026-capai-221-ribu'

$ws.Range("A17").Value = '5 Jurusan Kuliah yang Cocok Bekerja di Kapal Pesiar'
$ws.Range("B17").Value = 'Kampus | Selasa, 30 September 2025 12:14 WIB 5 Jurusan Kuliah yang Cocok Bekerja di Kapal Pesiar'
$ws.Range("D17").Value = 'https://edukasi.okezone.com/read/2025/09/30/65/3173445/5-jurusan-kuliah-yang-cocok-bekerja-di-kapal-pesiar'

$ws.Range("A18").Value = 'Gubernur Puji Damkar: 13 Jam Melawan Kobaran Api, 12 Kucing pun Ikut Selamat'
$ws.Range("B18").Value = 'Megapolitan | Selasa, 30 September 2025 12:08 WIB Gubernur Puji Damkar: 13 Jam Melawan Kobaran Api, 12 Kucing pun Ikut Selamat'
$ws.Range("D18").Value = 'https://news.okezone.com/read/2025/09/30/338/3173444/gubernur-puji-damkar-13-jam-melawan-kobaran-api-12-kucing-pun-ikut-selamat'

$ws.Range("A19").Value = 'Cegah Aktivitas Ilegal, RI-Australia Gelar Patroli di Perairan Perbatasan'
$ws.Range("B19").Value = 'Hot Issue | Selasa, 30 September 2025 12:07 WIB Cegah Aktivitas Ilegal, RI-Australia Gelar Patroli di Perairan Perbatasan'
$ws.Range("D19").Value = 'https://economy.okezone.com/read/2025/09/30/320/3173443/cegah-aktivitas-ilegal-ri-australia-gelar-patroli-di-perairan-perbatasan'

$ws.Range("A20").Value = 'PSSI Tak Risau soal Alokasi Anggaran Pemerintah untuk Timnas Indonesia U-23 di SEA Games 2025'
$ws.Range("B20").Value = 'Sepakbola Dunia | Selasa, 30 September 2025 12:05 WIB PSSI Tak Risau soal Alokasi Anggaran Pemerintah untuk Timnas Indonesia U-23 di SEA Games 2025'
$ws.Range("D20").Value = 'https://bola.okezone.com/read/2025/09/30/51/3173442/pssi-tak-risau-soal-alokasi-anggaran-pemerintah-untuk-timnas-indonesia-u-23-di-sea-games-2025'

$ws.Range("A21").Value = 'Tanggapan Keluarga Usai Tanggal Lahir Bing Slamet Jadi Hari Komedi Nasional'
$ws.Range("B21").Value = 'Hot Gossip | Selasa, 30 September 2025 12:01 WIB Tanggapan Keluarga Usai Tanggal Lahir Bing Slamet Jadi Hari Komedi Nasional'
$ws.Range("D21").Value = 'https://celebrity.okezone.com/read/2025/09/29/33/3173360/tanggapan-keluarga-usai-tanggal-lahir-bing-slamet-jadi-hari-komedi-nasional'

$ws.Range("A22").Value = 'Prediksi Chelsea vs Benfica di Liga Champions 2025-2026: Kembalinya Jose Mourinho ke Stamford Bridge'
$ws.Range("B22").Value = 'Liga Champion | Selasa, 30 September 2025 11:59 WIB Prediksi Chelsea vs Benfica di Liga Champions 2025-2026: Kembalinya Jose Mourinho ke Stamford Bridge'
$ws.Range("D22").Value = 'https://bola.okezone.com/read/2025/09/30/261/3173441/prediksi-chelsea-vs-benfica-di-liga-champions-2025-2026-kembalinya-jose-mourinho-ke-stamford-bridge'

$ws.Range("A23").Value = 'BGN Nonaktifkan Sementara 56 SPPG Imbas Adanya Kasus Keracunan MBG'
$ws.Range("B23").Value = 'Nasional | Selasa, 30 September 2025 11:54 WIB BGN Nonaktifkan Sementara 56 SPPG Imbas Adanya Kasus Keracunan MBG'
$ws.Range("D23").Value = 'https://news.okezone.com/read/2025/09/30/337/3173440/bgn-nonaktifkan-sementara-56-sppg-imbas-adanya-kasus-keracunan-mbg'

$ws.Range("A24").Value = 'Apa Saja Hadits tentang Palestina di Akhir Zaman?'
$ws.Range("B24").Value = 'Serba-serbi | Selasa, 30 September 2025 11:48 WIB Apa Saja Hadits tentang Palestina di Akhir Zaman?'
$ws.Range("D24").Value = 'https://muslim.okezone.com/read/2025/09/30/614/3173438/apa-saja-hadits-tentang-palestina-di-akhir-zaman'

$ws.Range("A25").Value = '5 Pemain Top Arab Saudi yang Wajib Diwaspadai Timnas Indonesia, Nomor 1  Rekan Setim Cristiano Ronaldo!'
$ws.Range("B25").Value = 'Sepakbola Dunia | Selasa, 30 September 2025 11:47 WIB 5 Pemain Top Arab Saudi yang Wajib Diwaspadai Timnas Indonesia, Nomor 1  Rekan Setim Cristiano Ronaldo!'
$ws.Range("D25").Value = 'https://bola.okezone.com/read/2025/09/30/51/3173437/5-pemain-top-arab-saudi-yang-wajib-diwaspadai-timnas-indonesia-nomor-1-rekan-setim-cristiano-ronaldo'

$ws.Range("A26").Value = '5 Pemain yang Berpotensi Dicoret Patrick Kluivert Jelang Timnas Indonesia Lawan Arab Saudi dan Irak, Nomor 1 Andalan Persib Bandung!'
$ws.Range("B26").Value = 'Sepakbola Dunia | Selasa, 30 September 2025 11:46 WIB 5 Pemain yang Berpotensi Dicoret Patrick Kluivert Jelang Timnas Indonesia Lawan Arab Saudi dan Irak, Nomor 1 Andalan Persib Bandung!'
$ws.Range("D26").Value = 'https://bola.okezone.com/read/2025/09/30/51/3173436/5-pemain-yang-berpotensi-dicoret-patrick-kluivert-jelang-timnas-indonesia-lawan-arab-saudi-dan-irak-nomor-1-andalan-persib-bandung'

$ws.Range("A27").Value = 'Ayah Arya Daru: Penyampaian Penyidik Polda Metro Belum Bisa Menenangkan Kami'
$ws.Range("B27").Value = 'Megapolitan | Selasa, 30 September 2025 11:44 WIB Ayah Arya Daru: Penyampaian Penyidik Polda Metro Belum Bisa Menenangkan Kami'
$ws.Range("D27").Value = 'https://news.okezone.com/read/2025/09/30/338/3173435/ayah-arya-daru-penyampaian-penyidik-polda-metro-belum-bisa-menenangkan-kami'

$ws.Range("A28").Value = 'Bantah Ada Tunggakan Subsidi BUMN 2024, Purbaya Minta Data Segera Dibereskan'
$ws.Range("B28").Value = 'Hot Issue | Selasa, 30 September 2025 11:44 WIB Bantah Ada Tunggakan Subsidi BUMN 2024, Purbaya Minta Data Segera Dibereskan'
$ws.Range("D28").Value = 'https://economy.okezone.com/read/2025/09/30/320/3173434/bantah-ada-tunggakan-subsidi-bumn-2024-purbaya-minta-data-segera-dibereskan'

$ws.Range("A29").Value = 'NOC Indonesia Ingatkan Jangan Terpancing Isu Sanksi FIFA untuk Malaysia'
$ws.Range("B29").Value = 'Sport Lain | Selasa, 30 September 2025 11:36 WIB NOC Indonesia Ingatkan Jangan Terpancing Isu Sanksi FIFA untuk Malaysia'
$ws.Range("D29").Value = 'https://sports.okezone.com/read/2025/09/30/43/3173433/noc-indonesia-ingatkan-jangan-terpancing-isu-sanksi-fifa-untuk-malaysia'

$ws.Range("A30").Value = 'KPK kembali Panggil Ilham Akbar Habibie terkait Kasus Pengadaan Iklan BJB'
$ws.Range("B30").Value = 'Nasional | Selasa, 30 September 2025 11:31 WIB KPK kembali Panggil Ilham Akbar Habibie terkait Kasus Pengadaan Iklan BJB'
$ws.Range("D30").Value = 'https://news.okezone.com/read/2025/09/30/337/3173432/kpk-kembali-panggil-ilham-akbar-habibie-terkait-kasus-pengadaan-iklan-bjb'

$ws.Range("A31").Value = 'Paul Pogba hingga PSSI-nya Turki Desak FIFA dan UEFA untuk Boikot Israel'
$ws.Range("B31").Value = 'Sepakbola Dunia | Selasa, 30 September 2025 11:23 WIB Paul Pogba hingga PSSI-nya Turki Desak FIFA dan UEFA untuk Boikot Israel'
$ws.Range("D31").Value = 'https://bola.okezone.com/read/2025/09/30/51/3173431/paul-pogba-hingga-pssi-nya-turki-desak-fifa-dan-uefa-untuk-boikot-israel'

$ws.Range("A32").Value = 'DPRD DKI Minta Ranpergub Kawasan Tanpa Rokok Segera Disiapkan'
$ws.Range("B32").Value = 'Megapolitan | Selasa, 30 September 2025 11:19 WIB DPRD DKI Minta Ranpergub Kawasan Tanpa Rokok Segera Disiapkan'
$ws.Range("D32").Value = 'https://news.okezone.com/read/2025/09/30/338/3173430/dprd-dki-minta-ranpergub-kawasan-tanpa-rokok-segera-disiapkan'

$ws.Range("A33").Value = 'Indonesia Sambut Baik Proposal Perdamaian Gaza Trump, Siap Kerja Sama dengan AS'
$ws.Range("B33").Value = 'International | Selasa, 30 September 2025 11:12 WIB Indonesia Sambut Baik Proposal Perdamaian Gaza Trump, Siap Kerja Sama dengan AS'
$ws.Range("D33").Value = 'https://news.okezone.com/read/2025/09/30/18/3173429/indonesia-sambut-baik-proposal-perdamaian-gaza-trump-siap-kerja-sama-dengan-as'

$ws.Range("A34").Value = 'Suzuki Siap Luncurkan Motor Listrik di Indonesia Tahun Depan, e-Access?'
$ws.Range("B34").Value = 'Autos | Selasa, 30 September 2025 11:10 WIB Suzuki Siap Luncurkan Motor Listrik di Indonesia Tahun Depan, e-Access?'
$ws.Range("D34").Value = 'https://ototekno.okezone.com/read/2025/09/30/15/3173428/suzuki-siap-luncurkan-motor-listrik-di-indonesia-tahun-depan-e-access'

$ws.Range("A35").Value = 'Polda Metro Musnahkan 1,14 Ton Narkoba Senilai Rp1 Triliun'
$ws.Range("B35").Value = 'Megapolitan | Selasa, 30 September 2025 11:09 WIB Polda Metro Musnahkan 1,14 Ton Narkoba Senilai Rp1 Triliun'
$ws.Range("D35").Value = 'https://news.okezone.com/read/2025/09/30/338/3173427/polda-metro-musnahkan-1-14-ton-narkoba-senilai-rp1-triliun'

$ws.Range("A36").Value = 'Wakasau Pimpin Inspeksi Pengadaan 6 Jet Tempur T-50i, 2 Unit Tiba di Indonesia November'
$ws.Range("B36").Value = 'Nasional | Selasa, 30 September 2025 11:06 WIB Wakasau Pimpin Inspeksi Pengadaan 6 Jet Tempur T-50i, 2 Unit Tiba di Indonesia November'
$ws.Range("D36").Value = 'https://news.okezone.com/read/2025/09/30/337/3173426/wakasau-pimpin-inspeksi-pengadaan-6-jet-tempur-t-50i-2-unit-tiba-di-indonesia-november'

$ws.Range("A37").Value = 'Update Ranking FIFA Timnas Indonesia jika Menang atas Arab Saudi dan Irak di Kualifikasi Piala Dunia 2026: Melesat Tajam!'
$ws.Range("B37").Value = 'Sepakbola Dunia | Selasa, 30 September 2025 11:04 WIB Update Ranking FIFA Timnas Indonesia jika Menang atas Arab Saudi dan Irak di Kualifikasi Piala Dunia 2026: Melesat Tajam!'
$ws.Range("D37").Value = 'https://bola.okezone.com/read/2025/09/30/51/3173425/update-ranking-fifa-timnas-indonesia-jika-menang-atas-arab-saudi-dan-irak-di-kualifikasi-piala-dunia-2026-melesat-tajam'

$ws.Range("A38").Value = '10 Perkara yang Membatalkan Sholat'
$ws.Range("B38").Value = 'Serba-serbi | Selasa, 30 September 2025 10:55 WIB 10 Perkara yang Membatalkan Sholat'
$ws.Range("D38").Value = 'https://muslim.okezone.com/read/2025/09/30/614/3173424/10-perkara-yang-membatalkan-sholat'

$ws.Range("A39").Value = 'Komedian Bedu Hadiri Sidang Cerai Perdana dengan Irma di PA Jaksel Hari Ini'
$ws.Range("B39").Value = 'Hot Gossip | Selasa, 30 September 2025 10:51 WIB Komedian Bedu Hadiri Sidang Cerai Perdana dengan Irma di PA Jaksel Hari Ini'
$ws.Range("D39").Value = 'https://celebrity.okezone.com/read/2025/09/30/33/3173423/komedian-bedu-hadiri-sidang-cerai-perdana-dengan-irma-di-pa-jaksel-hari-ini'

$ws.Range("A40").Value = 'Prediksi Kairat Almaty vs Real Madrid di Liga Champions 2025-2026: Si Anak Bawang Jadi Korban Pelampiasan?'
$ws.Range("B40").Value = 'Liga Champion | Selasa, 30 September 2025 10:51 WIB Prediksi Kairat Almaty vs Real Madrid di Liga Champions 2025-2026: Si Anak Bawang Jadi Korban Pelampiasan?'
$ws.Range("D40").Value = 'https://bola.okezone.com/read/2025/09/30/261/3173422/prediksi-kairat-almaty-vs-real-madrid-di-liga-champions-2025-2026-si-anak-bawang-jadi-korban-pelampiasan'

$ws.Range("A41").Value = 'Valentino Rossi dan Marc Marquez Kompak Barengan Ada di Jakarta Jelang MotoGP Mandalika 2025'
$ws.Range("B41").Value = 'MotoGP | Selasa, 30 September 2025 10:46 WIB Valentino Rossi dan Marc Marquez Kompak Barengan Ada di Jakarta Jelang MotoGP Mandalika 2025'
$ws.Range("D41").Value = 'https://sports.okezone.com/read/2025/09/30/38/3173420/valentino-rossi-dan-marc-marquez-kompak-barengan-ada-di-jakarta-jelang-motogp-mandalika-2025'

$ws.Range("A42").Value = 'Trump dan Netanyahu Sepakati Rencana Perdamaian Gaza, Ini Isinya'
$ws.Range("B42").Value = 'International | Selasa, 30 September 2025 10:43 WIB Trump dan Netanyahu Sepakati Rencana Perdamaian Gaza, Ini Isinya'
$ws.Range("D42").Value = 'https://news.okezone.com/read/2025/09/30/18/3173419/trump-dan-netanyahu-sepakati-rencana-perdamaian-gaza-ini-isinya'

$ws.Range("A43").Value = 'Butuh Obat Mujarab Kembalikan Penjualan Mobil Jadi 1 Juta Unit'
$ws.Range("B43").Value = 'Autos | Selasa, 30 September 2025 10:38 WIB Butuh Obat Mujarab Kembalikan Penjualan Mobil Jadi 1 Juta Unit'
$ws.Range("D43").Value = 'https://ototekno.okezone.com/read/2025/09/30/15/3173417/butuh-obat-mujarab-kembalikan-penjualan-mobil-jadi-1-juta-unit'

$ws.Range("A44").Value = 'Jadwal Siaran Langsung Timnas Indonesia di Babak 4 Kualifikasi Piala Dunia 2026 Zona Asia, Live di RCTI!'
$ws.Range("B44").Value = 'Sepakbola Dunia | Selasa, 30 September 2025 10:37 WIB Jadwal Siaran Langsung Timnas Indonesia di Babak 4 Kualifikasi Piala Dunia 2026 Zona Asia, Live di RCTI!'
$ws.Range("D44").Value = 'https://bola.okezone.com/read/2025/09/30/51/3173416/jadwal-siaran-langsung-timnas-indonesia-di-babak-4-kualifikasi-piala-dunia-2026-zona-asia-live-di-rcti'

$ws.Range("A45").Value = 'PLN Buka Lowongan Kerja 2025 untuk Lulusan D3-S2, Ini Cara Daftarnya'
$ws.Range("B45").Value = 'Smart Money | Selasa, 30 September 2025 10:34 WIB PLN Buka Lowongan Kerja 2025 untuk Lulusan D3-S2, Ini Cara Daftarnya'
$ws.Range("D45").Value = 'https://economy.okezone.com/read/2025/09/30/622/3173415/pln-buka-lowongan-kerja-2025-untuk-lulusan-d3-s2-ini-cara-daftarnya'

$ws.Range("A46").Value = 'Pramono: Api Cepat Membakar Rumah di Taman Sari Akibat Angin Kencang'
$ws.Range("B46").Value = 'Megapolitan | Selasa, 30 September 2025 10:32 WIB Pramono: Api Cepat Membakar Rumah di Taman Sari Akibat Angin Kencang'
$ws.Range("D46").Value = 'https://news.okezone.com/read/2025/09/30/338/3173414/pramono-api-cepat-membakar-rumah-di-taman-sari-akibat-angin-kencang'

$ws.Range("A47").Value = 'Herti Sastra Perjuangkan BPJS hingga Kesejahteraan Guru Ngaji di Deli Serdang Sumut'
$ws.Range("B47").Value = 'Nusantara | Selasa, 30 September 2025 10:27 WIB Herti Sastra Perjuangkan BPJS hingga Kesejahteraan Guru Ngaji di Deli Serdang Sumut'
$ws.Range("D47").Value = 'https://news.okezone.com/read/2025/09/30/340/3173413/herti-sastra-perjuangkan-bpjs-hingga-kesejahteraan-guru-ngaji-di-deli-serdang-sumut'

$ws.Range("A48").Value = 'Erick Thohir Jawab Tuduhan Intervensi FIFA hingga Bikin Malaysia Dihukum'
$ws.Range("B48").Value = 'Sepakbola Dunia | Selasa, 30 September 2025 10:21 WIB Erick Thohir Jawab Tuduhan Intervensi FIFA hingga Bikin Malaysia Dihukum'
$ws.Range("D48").Value = 'https://bola.okezone.com/read/2025/09/30/51/3173412/erick-thohir-jawab-tuduhan-intervensi-fifa-hingga-bikin-malaysia-dihukum'

$ws.Range("A49").Value = 'Pramono Janji Dampingi Pengurusan Surat HGB-SHM Korban Kebakaran Taman Sari'
$ws.Range("B49").Value = 'Megapolitan | Selasa, 30 September 2025 10:20 WIB Pramono Janji Dampingi Pengurusan Surat HGB-SHM Korban Kebakaran Taman Sari'
$ws.Range("D49").Value = 'https://news.okezone.com/read/2025/09/30/338/3173411/pramono-janji-dampingi-pengurusan-surat-hgb-shm-korban-kebakaran-taman-sari'

$ws.Range("A50").Value = 'KPPU Denda TikTok Rp15 Miliar Imbas Telat Lapor Akuisisi Tokopedia'
$ws.Range("B50").Value = 'Hot Issue | Selasa, 30 September 2025 10:11 WIB KPPU Denda TikTok Rp15 Miliar Imbas Telat Lapor Akuisisi Tokopedia'
$ws.Range("D50").Value = 'https://economy.okezone.com/read/2025/09/30/320/3173410/kppu-denda-tiktok-rp15-miliar-imbas-telat-lapor-akuisisi-tokopedia'

$ws.Range("A51").Value = 'Timnas Malaysia Kalah WO 0-3 dari Vietnam dan Nepal di Kualifikasi Piala Asia 2027 Setelah Mainkan Pemain Ilegal?'
$ws.Range("B51").Value = 'Sepakbola Dunia | Selasa, 30 September 2025 10:09 WIB Timnas Malaysia Kalah WO 0-3 dari Vietnam dan Nepal di Kualifikasi Piala Asia 2027 Setelah Mainkan Pemain Ilegal?'
$ws.Range("D51").Value = 'https://bola.okezone.com/read/2025/09/30/51/3173409/timnas-malaysia-kalah-wo-0-3-dari-vietnam-dan-nepal-di-kualifikasi-piala-asia-2027-setelah-mainkan-pemain-ilegal'

$ws.Range("A52").Value = 'Akses Produksi Pertanian Nagekeo Sempat Terputus, Kosmas Lawa Bagho Sigap Perbaiki Bersama Warga NTT'
$ws.Range("B52").Value = 'Nusantara | Selasa, 30 September 2025 10:06 WIB Akses Produksi Pertanian Nagekeo Sempat Terputus, Kosmas Lawa Bagho Sigap Perbaiki Bersama Warga NTT'
$ws.Range("D52").Value = 'https://news.okezone.com/read/2025/09/30/340/3173408/akses-produksi-pertanian-nagekeo-sempat-terputus-kosmas-lawa-bagho-sigap-perbaiki-bersama-warga-ntt'

$ws.Range("A53").Value = 'Nunung Akui Serakah soal Makanan Sebelum Idap Kanker'
$ws.Range("B53").Value = 'Hot Gossip | Selasa, 30 September 2025 10:03 WIB Nunung Akui Serakah soal Makanan Sebelum Idap Kanker'
$ws.Range("D53").Value = 'https://celebrity.okezone.com/read/2025/09/29/33/3173359/nunung-akui-serakah-soal-makanan-sebelum-idap-kanker'

$ws.Range("A54").Value = 'Hadir di Tengah Pertemuan Baim Wong dan Paula Verhoeven, Kimberly Ryder Bongkar Fakta Ini'
$ws.Range("B54").Value = 'Hot Gossip | Selasa, 30 September 2025 10:02 WIB Hadir di Tengah Pertemuan Baim Wong dan Paula Verhoeven, Kimberly Ryder Bongkar Fakta Ini'
$ws.Range("D54").Value = 'https://celebrity.okezone.com/read/2025/09/29/33/3173356/hadir-di-tengah-pertemuan-baim-wong-dan-paula-verhoeven-kimberly-ryder-bongkar-fakta-ini'

$ws.Range("A55").Value = '8 Santri Tertimbun Puing Mushola Ponpes di Sidoarjo Berhasil Dievakuasi'
$ws.Range("B55").Value = 'Nusantara | Selasa, 30 September 2025 10:00 WIB 8 Santri Tertimbun Puing Mushola Ponpes di Sidoarjo Berhasil Dievakuasi'
$ws.Range("D55").Value = 'https://news.okezone.com/read/2025/09/30/340/3173406/8-santri-tertimbun-puing-mushola-ponpes-di-sidoarjo-berhasil-dievakuasi'

$ws.Range("A56").Value = 'Sabrina Hapus Nama Deddy Corbuzier di Medsos, Ada Apa?'
$ws.Range("B56").Value = 'Life | Selasa, 30 September 2025 09:53 WIB Sabrina Hapus Nama Deddy Corbuzier di Medsos, Ada Apa?'
$ws.Range("D56").Value = 'https://women.okezone.com/read/2025/09/30/612/3173407/sabrina-hapus-nama-deddy-corbuzier-di-medsos-ada-apa'

$ws.Range("A57").Value = 'Purbaya Siap Buka-bukaan Subsidi BBM hingga Listrik di DPR'
$ws.Range("B57").Value = 'Hot Issue | Selasa, 30 September 2025 09:52 WIB Purbaya Siap Buka-bukaan Subsidi BBM hingga Listrik di DPR'
$ws.Range("D57").Value = 'https://economy.okezone.com/read/2025/09/30/320/3173405/purbaya-siap-buka-bukaan-subsidi-bbm-hingga-listrik-di-dpr'

$ws.Range("A58").Value = '5 Atlet Bulu Tangkis Dunia yang Suka Kuliner Khas Indonesia, Nomor 1 Doyan Bakso!'
$ws.Range("B58").Value = 'Netting | Selasa, 30 September 2025 09:44 WIB 5 Atlet Bulu Tangkis Dunia yang Suka Kuliner Khas Indonesia, Nomor 1 Doyan Bakso!'
$ws.Range("D58").Value = 'https://sports.okezone.com/read/2025/09/30/40/3173404/5-atlet-bulu-tangkis-dunia-yang-suka-kuliner-khas-indonesia-nomor-1-doyan-bakso'

$ws.Range("A59").Value = 'IHSG Dibuka Menguat ke 8.150, Dekati Rekor Tertinggi'
$ws.Range("B59").Value = 'Market Update | Selasa, 30 September 2025 09:40 WIB IHSG Dibuka Menguat ke 8.150, Dekati Rekor Tertinggi'
$ws.Range("D59").Value = 'https://economy.okezone.com/read/2025/09/30/278/3173402/ihsg-dibuka-menguat-ke-nbsp-8-150-dekati-rekor-tertinggi-nbsp'

$ws.Range("A60").Value = 'Waspada! Obesitas Tingkatkan Risiko Penyakit Jantung'
$ws.Range("B60").Value = 'Health | Selasa, 30 September 2025 09:37 WIB Waspada! Obesitas Tingkatkan Risiko Penyakit Jantung'
$ws.Range("D60").Value = 'https://women.okezone.com/read/2025/09/30/482/3173401/waspada-obesitas-tingkatkan-risiko-penyakit-jantung'

$ws.Range("A61").Value = 'Pramono Tinjau Lokasi Kebakaran di Taman Sari, Warga Butuh Popok dan Urus Surat Berharga'
$ws.Range("B61").Value = 'Megapolitan | Selasa, 30 September 2025 09:32 WIB Pramono Tinjau Lokasi Kebakaran di Taman Sari, Warga Butuh Popok dan Urus Surat Berharga'
$ws.Range("D61").Value = 'https://news.okezone.com/read/2025/09/30/338/3173400/pramono-tinjau-lokasi-kebakaran-di-taman-sari-warga-butuh-popok-dan-urus-surat-berharga'

$ws.Range("A62").Value = '4 Negara yang Pernah Terseret Skandal Naturalisasi, Nomor 1 Tetangga Indonesia!'
$ws.Range("B62").Value = 'Sepakbola Dunia | Selasa, 30 September 2025 09:22 WIB 4 Negara yang Pernah Terseret Skandal Naturalisasi, Nomor 1 Tetangga Indonesia!'
$ws.Range("D62").Value = 'https://bola.okezone.com/read/2025/09/30/51/3173398/4-negara-yang-pernah-terseret-skandal-naturalisasi-nomor-1-tetangga-indonesia'

$ws.Range("A63").Value = 'Pelaku Pembakaran Kantor Polisi di Kediri Ditangkap!'
$ws.Range("B63").Value = 'Nusantara | Selasa, 30 September 2025 09:16 WIB Pelaku Pembakaran Kantor Polisi di Kediri Ditangkap!'
$ws.Range("D63").Value = 'https://news.okezone.com/read/2025/09/30/340/3173397/pelaku-pembakaran-kantor-polisi-di-kediri-ditangkap'

$ws.Range("A64").Value = 'Kenapa Ikan Hiu Mengandung Merkuri yang Tinggi?'
$ws.Range("B64").Value = 'Food | Selasa, 30 September 2025 09:10 WIB Kenapa Ikan Hiu Mengandung Merkuri yang Tinggi?'
$ws.Range("D64").Value = 'https://women.okezone.com/read/2025/09/29/298/3173305/kenapa-ikan-hiu-mengandung-merkuri-yang-tinggi'

$ws.Range("A65").Value = 'Billy Syahputra dan Vika Kolesnaya Dikaruniai Anak Pertama Setelah Menikah'
$ws.Range("B65").Value = 'Hot Gossip | Selasa, 30 September 2025 09:01 WIB Billy Syahputra dan Vika Kolesnaya Dikaruniai Anak Pertama Setelah Menikah'
$ws.Range("D65").Value = 'https://celebrity.okezone.com/read/2025/09/29/33/3173354/billy-syahputra-dan-vika-kolesnaya-dikaruniai-anak-pertama-setelah-menikah-nbsp'

$ws.Range("A66").Value = '5 Fakta Mushola Roboh di Ponpes Sidoarjo, 1 Orang Tewas dan 79 Terluka'
$ws.Range("B66").Value = 'Nasional | Selasa, 30 September 2025 08:59 WIB 5 Fakta Mushola Roboh di Ponpes Sidoarjo, 1 Orang Tewas dan 79 Terluka'
$ws.Range("D66").Value = 'https://news.okezone.com/read/2025/09/30/337/3173395/5-fakta-mushola-roboh-di-ponpes-sidoarjo-1-orang-tewas-dan-79-terluka'

$ws.Range("A67").Value = 'Harga Emas Antam Hari Ini Cetak Rekor Lagi, Tembus Rp2,2 Juta per Gram!'
$ws.Range("B67").Value = 'Hot Issue | Selasa, 30 September 2025 08:57 WIB Harga Emas Antam Hari Ini Cetak Rekor Lagi, Tembus Rp2,2 Juta per Gram!'
$ws.Range("D67").Value = 'https://economy.okezone.com/read/2025/09/30/320/3173394/harga-emas-antam-hari-ini-cetak-rekor-lagi-tembus-rp2-2-juta-per-gram'

$ws.Range("A68").Value = 'Baru Jadi Juara Dunia MotoGP 2025, Marc Marquez Langsung Ditantang Taklukkan Sirkuit Mandalika!'
$ws.Range("B68").Value = 'MotoGP | Selasa, 30 September 2025 08:50 WIB Baru Jadi Juara Dunia MotoGP 2025, Marc Marquez Langsung Ditantang Taklukkan Sirkuit Mandalika!'
$ws.Range("D68").Value = 'https://sports.okezone.com/read/2025/09/30/38/3173393/baru-jadi-juara-dunia-motogp-2025-marc-marquez-langsung-ditantang-taklukkan-sirkuit-mandalika'

$ws.Range("A69").Value = 'Sinopsis Series Vision+ Cinta di Balik Awan Episode 2, Hanya di RCTI'
$ws.Range("B69").Value = 'TV Scoop | Selasa, 30 September 2025 08:30 WIB Sinopsis Series Vision+ Cinta di Balik Awan Episode 2, Hanya di RCTI'
$ws.Range("D69").Value = 'https://celebrity.okezone.com/read/2025/09/30/598/3173399/sinopsis-series-vision-cinta-di-balik-awan-episode-2-hanya-di-rcti'

$ws.Range("A70").Value = 'Raperda Kawasan Tanpa Rokok Dikebut DPRD, Ini Kata Orang Dekat Pramono'
$ws.Range("B70").Value = 'Megapolitan | Selasa, 30 September 2025 08:10 WIB Raperda Kawasan Tanpa Rokok Dikebut DPRD, Ini Kata Orang Dekat Pramono'
$ws.Range("D70").Value = 'https://news.okezone.com/read/2025/09/30/338/3173392/raperda-kawasan-tanpa-rokok-dikebut-dprd-ini-kata-orang-dekat-pramono'

$ws.Range("A71").Value = 'Potret Ruben Onsu Cium Hajar Aswad'
$ws.Range("B71").Value = 'Life | Selasa, 30 September 2025 08:10 WIB Potret Ruben Onsu Cium Hajar Aswad'
$ws.Range("D71").Value = 'https://women.okezone.com/read/2025/09/29/612/3173301/potret-ruben-onsu-cium-hajar-aswad'

$ws.Range("A72").Value = 'PSSI: Gerald Vanenburg Tetap Pelatih Timnas Indonesia U-23, Indra Sjafri Hanya untuk SEA Games 2025'
$ws.Range("B72").Value = 'Sepakbola Dunia | Selasa, 30 September 2025 08:06 WIB PSSI: Gerald Vanenburg Tetap Pelatih Timnas Indonesia U-23, Indra Sjafri Hanya untuk SEA Games 2025'
$ws.Range("D72").Value = 'https://bola.okezone.com/read/2025/09/30/51/3173391/pssi-gerald-vanenburg-tetap-pelatih-timnas-indonesia-u-23-indra-sjafri-hanya-untuk-sea-games-2025'

$ws.Range("A73").Value = 'Zulhas Sebut Tak Perlu Semua Masalah Dibawa ke Prabowo'
$ws.Range("B73").Value = 'Hot Issue | Selasa, 30 September 2025 08:02 WIB Zulhas Sebut Tak Perlu Semua Masalah Dibawa ke Prabowo'
$ws.Range("D73").Value = 'https://economy.okezone.com/read/2025/09/30/320/3173390/zulhas-sebut-tak-perlu-semua-masalah-dibawa-ke-prabowo'

$ws.Range("A74").Value = '3 Pembalap MotoGP yang Pernah Menang di Sirkuit Mandalika, Nomor 1 Namanya Terukir dalam Sejarah!'
$ws.Range("B74").Value = 'MotoGP | Selasa, 30 September 2025 07:55 WIB 3 Pembalap MotoGP yang Pernah Menang di Sirkuit Mandalika, Nomor 1 Namanya Terukir dalam Sejarah!'
$ws.Range("D74").Value = 'https://sports.okezone.com/read/2025/09/30/38/3173389/3-pembalap-motogp-yang-pernah-menang-di-sirkuit-mandalika-nomor-1-namanya-terukir-dalam-sejarah'

$ws.Range("A75").Value = 'Patrick Kluivert Panggil Cyrus Margono untuk Laga Timnas Indonesia vs Arab Saudi dan Irak karena Emil Audero Cedera?'
$ws.Range("B75").Value = 'Sepakbola Dunia | Selasa, 30 September 2025 07:38 WIB Patrick Kluivert Panggil Cyrus Margono untuk Laga Timnas Indonesia vs Arab Saudi dan Irak karena Emil Audero Cedera?'
$ws.Range("D75").Value = 'https://bola.okezone.com/read/2025/09/30/51/3173388/patrick-kluivert-panggil-cyrus-margono-untuk-laga-timnas-indonesia-vs-arab-saudi-dan-irak-karena-emil-audero-cedera'

# Rows where relevan/keywords_found changed from True/keyword to False/empty
$ws.Range("E34").Value = $false
$ws.Range("F34").ClearContents()
$ws.Range("E44").Value = $false
$ws.Range("F44").ClearContents()
$ws.Range("E53").Value = $false
$ws.Range("F53").ClearContents()
$ws.Range("E74").Value = $false
$ws.Range("F74").ClearContents()

Write-Host "Edit complete"
